$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 679.6585
$ws.Range("J17").Value = 577.75757
$ws.Range("L17").Value = 1733.27271
$ws.Range("N17").Value = -2069.27271

$ws.Range("H53").Value = 451.07144
$ws.Range("I53").Value = 341.6154
$ws.Range("J53").Value = 545.93335
$ws.Range("K53").Value = 341.6154
$ws.Range("L53").Value = 545.93335
$ws.Range("M53").Value = 295.3846
$ws.Range("N53").Value = -1819.93335

$ws.Range("H62").Value = 3101.3635
$ws.Range("I62").Value = 1761.6666
$ws.Range("J62").Value = 3603.75
$ws.Range("K62").Value = 1761.6666
$ws.Range("L62").Value = 3603.75
$ws.Range("M62").Value = -1137.6666
$ws.Range("N62").Value = -4851.75

$ws.Range("H65").Value = 3101.3635
$ws.Range("I65").Value = 1761.6666
$ws.Range("J65").Value = 3603.75
$ws.Range("K65").Value = 8808.333000000001
$ws.Range("L65").Value = 18018.75
$ws.Range("M65").Value = -5688.333000000001
$ws.Range("N65").Value = -24258.75

$ws.Range("H107").Value = 1933.04
$ws.Range("I107").Value = 1461.5294
$ws.Range("J107").Value = 2935
$ws.Range("K107").Value = 1461.5294
$ws.Range("L107").Value = 2935
$ws.Range("M107").Value = 458.4706000000001
$ws.Range("N107").Value = -6775

$ws.Range("H125").Value = 2446.389
$ws.Range("J125").Value = 2957.818
$ws.Range("L125").Value = 26620.362
$ws.Range("N125").Value = -31540.362

$ws.Range("H129").Value = 873.8823
$ws.Range("J129").Value = 958.9589
$ws.Range("L129").Value = 2876.8767
$ws.Range("N129").Value = -12876.8767

$ws.Range("H132").Value = 24636432
$ws.Range("I132").Value = 30304694
$ws.Range("J132").Value = 1254851.5
$ws.Range("K132").Value = 90914082
$ws.Range("L132").Value = 3764554.5
$ws.Range("M132").Value = -90911552
$ws.Range("N132").Value = -3769614.5

$ws.Range("H137").Value = 2260.68
$ws.Range("I137").Value = 1180.6875
$ws.Range("K137").Value = 3542.0625
$ws.Range("M137").Value = -992.0625

$ws.Range("H138").Value = 2898
$ws.Range("I138").Value = 884.35
$ws.Range("J138").Value = 3401.4126
$ws.Range("K138").Value = 2653.05
$ws.Range("L138").Value = 10204.2378
$ws.Range("M138").Value = 2486.95
$ws.Range("N138").Value = -20484.2378

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3875.0857
$ws.Range("I32").Value = 3509.9841
$ws.Range("K32").Value = 3509.9841
$ws.Range("M32").Value = -3222.9841

$ws.Range("H61").Value = 767.3
$ws.Range("I61").Value = 559.2
$ws.Range("J61").Value = 1391.6
$ws.Range("K61").Value = 559.2
$ws.Range("L61").Value = 1391.6
$ws.Range("M61").Value = -347.2
$ws.Range("N61").Value = -1815.6

$ws.Range("H74").Value = 2454.1965
$ws.Range("I74").Value = 2457.9111
$ws.Range("J74").Value = 2439
$ws.Range("K74").Value = 2457.9111
$ws.Range("L74").Value = 2439
$ws.Range("M74").Value = -1583.9111
$ws.Range("N74").Value = -4187

$ws.Range("H77").Value = 2454.1965
$ws.Range("I77").Value = 2457.9111
$ws.Range("J77").Value = 2439
$ws.Range("K77").Value = 12289.5555
$ws.Range("L77").Value = 12195
$ws.Range("M77").Value = -7921.555499999999
$ws.Range("N77").Value = -20931

$ws.Range("H132").Value = 1685.0934
$ws.Range("I132").Value = 1147.9454
$ws.Range("J132").Value = 3162.25
$ws.Range("K132").Value = 3443.8362
$ws.Range("L132").Value = 9486.75
$ws.Range("M132").Value = -913.8362000000002
$ws.Range("N132").Value = -14546.75

$ws.Range("H136").Value = 767.3
$ws.Range("I136").Value = 559.2
$ws.Range("J136").Value = 1391.6
$ws.Range("K136").Value = 1677.6
$ws.Range("L136").Value = 4174.799999999999
$ws.Range("M136").Value = 872.3999999999999
$ws.Range("N136").Value = -9274.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1891.7462
$ws.Range("I134").Value = 1078.82
$ws.Range("J134").Value = 4282.706
$ws.Range("K134").Value = 3236.46
$ws.Range("L134").Value = 12848.118
$ws.Range("M134").Value = -701.46
$ws.Range("N134").Value = -17918.118

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8199211
$ws.Range("I31").Value = 1289.4517
$ws.Range("J31").Value = 16670397
$ws.Range("K31").Value = 1289.4517
$ws.Range("L31").Value = 16670397
$ws.Range("M31").Value = -994.4517000000001
$ws.Range("N31").Value = -16670987

$ws.Range("H34").Value = 8199211
$ws.Range("I34").Value = 1289.4517
$ws.Range("J34").Value = 16670397
$ws.Range("K34").Value = 1289.4517
$ws.Range("L34").Value = 16670397
$ws.Range("M34").Value = -1087.4517
$ws.Range("N34").Value = -16670801

$ws.Range("H58").Value = 1249.37
$ws.Range("I58").Value = 1273.0596
$ws.Range("J58").Value = 1125
$ws.Range("K58").Value = 1273.0596
$ws.Range("L58").Value = 1125
$ws.Range("M58").Value = -1070.0596
$ws.Range("N58").Value = -1531

$ws.Range("H134").Value = 1827.5526
$ws.Range("I134").Value = 943.7
$ws.Range("J134").Value = 2809.611
$ws.Range("K134").Value = 2831.1
$ws.Range("L134").Value = 8428.832999999999
$ws.Range("M134").Value = -296.1000000000004
$ws.Range("N134").Value = -13498.833

$ws.Range("H136").Value = 1249.37
$ws.Range("I136").Value = 1273.0596
$ws.Range("J136").Value = 1125
$ws.Range("K136").Value = 3819.1788
$ws.Range("L136").Value = 3375
$ws.Range("M136").Value = -1269.1788
$ws.Range("N136").Value = -8475

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 444.83334
$ws.Range("I97").Value = 262.4
$ws.Range("J97").Value = 575.1429000000001
$ws.Range("K97").Value = 787.1999999999999
$ws.Range("L97").Value = 1725.4287
$ws.Range("M97").Value = -291.1999999999999
$ws.Range("N97").Value = -2717.4287

$ws.Range("H131").Value = 869.62823
$ws.Range("I131").Value = 756.6667
$ws.Range("J131").Value = 879.0417
$ws.Range("K131").Value = 2270.0001
$ws.Range("L131").Value = 2637.1251
$ws.Range("M131").Value = 2769.9999
$ws.Range("N131").Value = -12717.1251

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 11595.235
$ws.Range("J123").Value = 11595.235
$ws.Range("L123").Value = 11595.235
$ws.Range("N123").Value = -16495.235

$ws.Range("H132").Value = 2732.3635
$ws.Range("I132").Value = 1782.6364
$ws.Range("J132").Value = 4631.8184
$ws.Range("K132").Value = 5347.9092
$ws.Range("L132").Value = 13895.4552
$ws.Range("M132").Value = -2817.9092
$ws.Range("N132").Value = -18955.4552

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1724.2354
$ws.Range("I46").Value = 1538.75
$ws.Range("J46").Value = 1889.1111
$ws.Range("K46").Value = 1538.75
$ws.Range("L46").Value = 1889.1111
$ws.Range("M46").Value = -1350.75
$ws.Range("N46").Value = -2265.1111

$ws.Range("H100").Value = 1800.1818
$ws.Range("I100").Value = 1500.25
$ws.Range("K100").Value = 1500.25
$ws.Range("M100").Value = -959.25

$ws.Range("H132").Value = 10786.763
$ws.Range("I132").Value = 11057.615
$ws.Range("K132").Value = 33172.845
$ws.Range("M132").Value = -30642.845

$ws.Range("H136").Value = 1971.2963
$ws.Range("I136").Value = 1065.119
$ws.Range("J136").Value = 5142.9165
$ws.Range("K136").Value = 3195.357
$ws.Range("L136").Value = 15428.7495
$ws.Range("M136").Value = -645.357
$ws.Range("N136").Value = -20528.7495

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 60237.5
$ws.Range("J46").Value = 60237.5
$ws.Range("L46").Value = 60237.5
$ws.Range("N46").Value = -60699.5

$ws.Range("H132").Value = 7093499
$ws.Range("I132").Value = 509.9189
$ws.Range("J132").Value = 33337558
$ws.Range("K132").Value = 1529.7567
$ws.Range("L132").Value = 100012674
$ws.Range("M132").Value = 1000.2433
$ws.Range("N132").Value = -100017734

$ws.Range("H134").Value = 60237.5
$ws.Range("J134").Value = 60237.5
$ws.Range("L134").Value = 180712.5
$ws.Range("N134").Value = -185782.5

$ws.Range("H136").Value = 1578.4133
$ws.Range("I136").Value = 463.29092
$ws.Range("J136").Value = 4645
$ws.Range("K136").Value = 1389.87276
$ws.Range("L136").Value = 13935
$ws.Range("M136").Value = 1160.12724
$ws.Range("N136").Value = -19035
